$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 13.30571766666667
$ws.Range("H2").Value = 39.917153
$ws.Range("I2").Value = 0.007643519924167935
$ws.Range("J2").Value = 0.007643519924167933
$ws.Range("M2").Value = 2.680851666666667
$ws.Range("N2").Value = 8.042555
$ws.Range("O2").Value = 0.1074910720871699
$ws.Range("P2").Value = 0.1074910720871699
$ws.Range("Q2").Value = 35.67065538287945
$ws.Range("R2").Value = 321.035898445915
$ws.Range("S2").Value = 0.000821610151168455
$ws.Range("T2").Value = 0.0008216101511684546

# Row 3
$ws.Range("G3").Value = 13.30571766666667
$ws.Range("H3").Value = 39.917153
$ws.Range("I3").Value = 0.007643519924167935
$ws.Range("J3").Value = 0.007643519924167933
$ws.Range("M3").Value = 18.51427066666667
$ws.Range("O3").Value = 0.7423457357290222
$ws.Range("P3").Value = 0.7423457357290222
$ws.Range("Q3").Value = 246.3456582949151
$ws.Range("R3").Value = 2217.110924654236
$ws.Range("S3").Value = 0.005674134421665885
$ws.Range("T3").Value = 0.005674134421665884

# Row 4
$ws.Range("G4").Value = 13.30571766666667
$ws.Range("H4").Value = 39.917153
$ws.Range("I4").Value = 0.007643519924167935
$ws.Range("J4").Value = 0.007643519924167933
$ws.Range("M4").Value = 3.745104
$ws.Range("N4").Value = 11.235312
$ws.Range("O4").Value = 0.1501631921838079
$ws.Range("P4").Value = 0.1501631921838079
$ws.Range("Q4").Value = 49.831296456304
$ws.Range("R4").Value = 448.481668106736
$ws.Range("S4").Value = 0.001147775351333594
$ws.Range("T4").Value = 0.001147775351333594

# Row 5
$ws.Range("I5").Value = 0.9448263940026712
$ws.Range("J5").Value = 0.9448263940026712
$ws.Range("M5").Value = 2.680851666666667
$ws.Range("N5").Value = 8.042555
$ws.Range("O5").Value = 0.1074910720871699
$ws.Range("P5").Value = 0.1074910720871699
$ws.Range("Q5").Value = 4409.300561977248
$ws.Range("R5").Value = 39683.70505779523
$ws.Range("S5").Value = 0.1015604020276019
$ws.Range("T5").Value = 0.1015604020276019

# Row 6
$ws.Range("I6").Value = 0.9448263940026712
$ws.Range("J6").Value = 0.9448263940026712
$ws.Range("M6").Value = 18.51427066666667
$ws.Range("O6").Value = 0.7423457357290222
$ws.Range("P6").Value = 0.7423457357290222
$ws.Range("S6").Value = 0.7013878445921119
$ws.Range("T6").Value = 0.7013878445921119

# Row 7
$ws.Range("I7").Value = 0.9448263940026712
$ws.Range("J7").Value = 0.9448263940026712
$ws.Range("M7").Value = 3.745104
$ws.Range("N7").Value = 11.235312
$ws.Range("O7").Value = 0.1501631921838079
$ws.Range("P7").Value = 0.1501631921838079
$ws.Range("Q7").Value = 6159.717591684447
$ws.Range("R7").Value = 55437.45832516003
$ws.Range("S7").Value = 0.1418781473829573
$ws.Range("T7").Value = 0.1418781473829573

# Row 8
$ws.Range("G8").Value = 82.73961633333333
$ws.Range("H8").Value = 248.218849
$ws.Range("I8").Value = 0.04753008607316088
$ws.Range("J8").Value = 0.04753008607316087
$ws.Range("M8").Value = 2.680851666666667
$ws.Range("N8").Value = 8.042555
$ws.Range("O8").Value = 0.1074910720871699
$ws.Range("P8").Value = 0.1074910720871699
$ws.Range("Q8").Value = 221.8126383465772
$ws.Range("R8").Value = 1996.313745119195
$ws.Range("S8").Value = 0.005109059908399527
$ws.Range("T8").Value = 0.005109059908399525

# Row 9
$ws.Range("G9").Value = 82.73961633333333
$ws.Range("H9").Value = 248.218849
$ws.Range("I9").Value = 0.04753008607316088
$ws.Range("J9").Value = 0.04753008607316087
$ws.Range("M9").Value = 18.51427066666667
$ws.Range("O9").Value = 0.7423457357290222
$ws.Range("P9").Value = 0.7423457357290222
$ws.Range("Q9").Value = 1531.863651651488
$ws.Range("R9").Value = 13786.77286486339
$ws.Range("S9").Value = 0.03528375671524436
$ws.Range("T9").Value = 0.03528375671524436

# Row 10
$ws.Range("G10").Value = 82.73961633333333
$ws.Range("H10").Value = 248.218849
$ws.Range("I10").Value = 0.04753008607316088
$ws.Range("J10").Value = 0.04753008607316087
$ws.Range("M10").Value = 3.745104
$ws.Range("N10").Value = 11.235312
$ws.Range("O10").Value = 0.1501631921838079
$ws.Range("P10").Value = 0.1501631921838079
$ws.Range("Q10").Value = 309.868468088432
$ws.Range("R10").Value = 2788.816212795888
$ws.Range("S10").Value = 0.007137269449516988
$ws.Range("T10").Value = 0.007137269449516987
